$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 19
$ws.Range("C4").Value = 15.9
$ws.Range("C5").Value = 22.5

$ws.Range("C5").Select()
